# This script reproduces the edit described by the diff:
#  - 17 new shared strings are added (used across the "Oplossingen" and
#    "Handelingen" sheets)
#  - 5 new rows are appended to the "Oplossingen" sheet (sheet2)
#  - 14 new rows are appended to the "Handelingen" sheet (sheet3)
#  - view/selection state is updated on all three sheets
#
# The order in which brand-new string VALUES are first assigned to cells
# below is deliberately chosen so that Excel's shared-string table ends up
# allocating indices 147-163 in the exact sequence required by the target
# workbook (new values are appended to the shared-string table in
# first-use order).

$wb = $excel.ActiveWorkbook

$wsInc = $wb.Worksheets.Item("Incidenten")
$wsOpl = $wb.Worksheets.Item("Oplossingen")
$wsHan = $wb.Worksheets.Item("Handelingen")

# ---------------------------------------------------------------------
# Oplossingen (sheet2): rows 31-32 only reuse already existing strings.
# ---------------------------------------------------------------------
$wsOpl.Cells.Item(31,1).Value = 30
$wsOpl.Cells.Item(31,2).Value = 12
$wsOpl.Cells.Item(31,3).Value = "Volgens protocol: probleem telefonisch melden aan team Multimedia"
$wsOpl.Cells.Item(31,4).Value = "Wachter op terugkoppeling van MMT, pas dan andere opties overwegen"

$wsOpl.Cells.Item(32,1).Value = 31
$wsOpl.Cells.Item(32,2).Value = 12
$wsOpl.Cells.Item(32,3).Value = "Meld probleem rechtsreeks aan Jacot"
$wsOpl.Cells.Item(32,4).Value = "Alleen wanneer MMT niet beschikbaar is"

# Row 33: D33 (string 147) is populated before C33 (string 148), matching
# the first-use order required for index assignment.
$wsOpl.Cells.Item(33,1).Value = 32
$wsOpl.Cells.Item(33,2).Value = 12
$wsOpl.Cells.Item(33,4).Value = "Bij het noodscenario wordt een opname automatisch gestart maar MVI er kan niet geklikt worden op de agendapunten, spreekttijd en stemmingen werken niet"
$wsOpl.Cells.Item(33,3).Value = "ALLEEN IN OVERLEG MET JACOT: Schakkel het noodscenario in"

# ---------------------------------------------------------------------
# Handelingen (sheet3): rows 86-88 reuse existing strings, row 89 adds two
# new strings (149, 150).
# ---------------------------------------------------------------------
$wsHan.Cells.Item(86,1).Value = 85
$wsHan.Cells.Item(86,2).Value = 30
$wsHan.Cells.Item(86,3).Value = "Bel een van de medewerkers van MMT en meld het probleem "
$wsHan.Cells.Item(86,4).Value = "Secretaris"

$wsHan.Cells.Item(87,1).Value = 86
$wsHan.Cells.Item(87,2).Value = 30
$wsHan.Cells.Item(87,3).Value = "MMT neemt contact op met Jacot en koppelt terug"
$wsHan.Cells.Item(87,4).Value = "MMT"
$wsHan.Cells.Item(87,3).Font.Name = "Bolder"
$wsHan.Cells.Item(87,3).Font.Size = 10

$wsHan.Cells.Item(88,1).Value = 87
$wsHan.Cells.Item(88,2).Value = 31
$wsHan.Cells.Item(88,3).Value = "Beld het noodnunmer van Jacot en meld het probleem"
$wsHan.Cells.Item(88,4).Value = "Secretaris"
$wsHan.Cells.Item(88,3).Font.Name = "Bolder"
$wsHan.Cells.Item(88,3).Font.Size = 10

$wsHan.Cells.Item(89,1).Value = 88
$wsHan.Cells.Item(89,2).Value = 32
$wsHan.Cells.Item(89,3).Value = "In Crestron, druk op de Home knop"
$wsHan.Cells.Item(89,4).Value = "Secretaris"
$wsHan.Cells.Item(89,5).Value = "Home knop Crestron.png"
$wsHan.Cells.Item(89,3).Font.Name = "Bolder"
$wsHan.Cells.Item(89,3).Font.Size = 10

$wsHan.Cells.Item(90,1).Value = 89
$wsHan.Cells.Item(90,2).Value = 32
$wsHan.Cells.Item(90,3).Value = "Druk drie seconden op de knop 'NOODFunctie camera uitschakelen'"
$wsHan.Cells.Item(90,4).Value = "Secretaris"
$wsHan.Cells.Item(90,3).Font.Name = "Bolder"
$wsHan.Cells.Item(90,3).Font.Size = 10

$wsHan.Cells.Item(91,1).Value = 90
$wsHan.Cells.Item(91,2).Value = 32
$wsHan.Cells.Item(91,3).Value = "Pak de twee handmicrofoons en doe deze aan via de knop aan de de achterkant van de microfoons"
$wsHan.Cells.Item(91,4).Value = "Secretaris"
$wsHan.Cells.Item(91,3).Font.Name = "Bolder"
$wsHan.Cells.Item(91,3).Font.Size = 10

$wsHan.Cells.Item(92,1).Value = 91
$wsHan.Cells.Item(92,2).Value = 32
$wsHan.Cells.Item(92,3).Value = "Laat twee bodes met de handmicrofoons rondlopen"
$wsHan.Cells.Item(92,4).Value = "Secretaris"
$wsHan.Cells.Item(92,3).Font.Name = "Bolder"
$wsHan.Cells.Item(92,3).Font.Size = 10

$wsHan.Cells.Item(93,1).Value = 92
$wsHan.Cells.Item(93,2).Value = 32
$wsHan.Cells.Item(93,3).Value = "Laat Jacot controleren of de opname loopt"
$wsHan.Cells.Item(93,4).Value = "Secretaris"
$wsHan.Cells.Item(93,3).Font.Name = "Bolder"
$wsHan.Cells.Item(93,3).Font.Size = 10

# ---------------------------------------------------------------------
# Back to Oplossingen: row 34 adds strings 155, 156; row 35 adds 160, 161.
# ---------------------------------------------------------------------
$wsOpl.Cells.Item(34,1).Value = 33
$wsOpl.Cells.Item(34,2).Value = 13
$wsOpl.Cells.Item(34,3).Value = "Verlaat direct de ruimte"
$wsOpl.Cells.Item(34,4).Value = "Het systeem zorgt er automatisch voor dat de vergadering wordt afgesloten. De opname tot dan toe wordt gered"

# ---------------------------------------------------------------------
# Handelingen rows 94-97 (existing strings), row 98 adds 157, 158, 159.
# ---------------------------------------------------------------------
$wsHan.Cells.Item(94,1).Value = 93
$wsHan.Cells.Item(94,2).Value = 32
$wsHan.Cells.Item(94,3).Value = "Via iBabs, publiek informeren dat de webcast later wordt geplaatst"
$wsHan.Cells.Item(94,4).Value = "Secretaris"

$wsHan.Cells.Item(95,1).Value = 94
$wsHan.Cells.Item(95,2).Value = 32
$wsHan.Cells.Item(95,3).Value = "Na de vergadering:"
$wsHan.Cells.Item(95,4).Value = "Secretaris"
$wsHan.Cells.Item(95,6).Value = "Bestanden beheren van vergaderingen.pdf"
$wsHan.Cells.Item(95,3).WrapText = $true
$wsHan.Rows.Item(95).RowHeight = 30

$wsHan.Cells.Item(96,1).Value = 95
$wsHan.Cells.Item(96,2).Value = 32
$wsHan.Cells.Item(96,3).Value = "Na de vergadering: "
$wsHan.Cells.Item(96,4).Value = "Secretaris"
$wsHan.Cells.Item(96,3).WrapText = $true

$wsHan.Cells.Item(97,1).Value = 96
$wsHan.Cells.Item(97,2).Value = 32
$wsHan.Cells.Item(97,3).Value = "Na de vergadering"
$wsHan.Cells.Item(97,4).Value = "Secretaris"
$wsHan.Cells.Item(97,3).WrapText = $true

$wsHan.Cells.Item(98,1).Value = 97
$wsHan.Cells.Item(98,2).Value = 33
$wsHan.Cells.Item(98,3).Value = "Als je dit beeld ziet op Crestron, verlaat de ruimte per direct"
$wsHan.Cells.Item(98,4).Value = "Allen"
$wsHan.Cells.Item(98,5).Value = "Brandmelding.png"
$wsHan.Cells.Item(98,3).WrapText = $true

# ---------------------------------------------------------------------
# Back to Oplossingen: row 35 adds strings 160, 161.
# ---------------------------------------------------------------------
$wsOpl.Cells.Item(35,1).Value = 34
$wsOpl.Cells.Item(35,2).Value = 14
$wsOpl.Cells.Item(35,3).Value = "Wacht max 5 minuten, indien geen stroom sluit de vergadering af"
$wsOpl.Cells.Item(35,4).Value = "Er is stroom voor maximaal 15 minuten, als de vergadering niet op tijd afgesloten dan gaat de opname tot dan toe verloren"

# ---------------------------------------------------------------------
# Handelingen row 99 adds strings 162, 163 (164th/last new string).
# ---------------------------------------------------------------------
$wsHan.Cells.Item(99,1).Value = 98
$wsHan.Cells.Item(99,2).Value = 34
$wsHan.Cells.Item(99,3).Value = "Bij een stroomstoring waarbij het systeem overgaat naar noodstroom, wordt deze melding zichtbaar in Crestron"
$wsHan.Cells.Item(99,4).Value = "MMT/Secretaris"
$wsHan.Cells.Item(99,5).Value = "Stroomstoring.png"
$wsHan.Cells.Item(99,3).WrapText = $true

# ---------------------------------------------------------------------
# View/selection state updates to match the final saved workbook.
# ---------------------------------------------------------------------
$wsInc.Activate()
$wsInc.Range("B15").Select()

$wsOpl.Activate()
$excel.ActiveWindow.ScrollRow = 16
$wsOpl.Range("C36").Select()

$wsHan.Activate()
$excel.ActiveWindow.ScrollRow = 82
$wsHan.Range("E93").Select()
